$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixing Iran FC Team Icons: rows 31-46 (S=Farsi name, T=English name).
# Cell write order matters because it determines the order new strings are
# appended to the shared-strings table (mirrors how the original edit was made).

# Row 31: Esteghlal -> Perspolis (reuses an existing shared string)
$ws.Range("S31").Value = "پرسپولیس تهران"
$ws.Range("T31").Value = "Perspolis"

# Row 32: Perspolis -> Sepahan (reuses an existing shared string)
$ws.Range("S32").Value = "سپاهان"
$ws.Range("T32").Value = "Sepahan"

# Row 33: Sepahan -> Esteghlal (reuses an existing shared string)
$ws.Range("S33").Value = "استقلال تهران"
$ws.Range("T33").Value = "Esteghlal"

# Row 34: Foolad -> Padideh Shahr Khodro (brand-new strings)
$ws.Range("T34").Value = "Padideh Shahr Khodro"
$ws.Range("S34").Value = "پدیده شهر خودرو"

# Row 35: Saipa -> Tractorsazi (reuses an existing shared string)
$ws.Range("S35").Value = "تراکتور سازی"
$ws.Range("T35").Value = "Tractorsazi"

# Row 36: Rah Ahan Tehran -> Zob Ahan (reuses an existing shared string)
$ws.Range("S36").Value = "ذوب آهن"
$ws.Range("T36").Value = "Zob Ahan"

# Row 37: Zob Ahan -> Saipa (reuses an existing shared string)
$ws.Range("S37").Value = "سایپا"
$ws.Range("T37").Value = "Saipa"

# Row 38: Tractorsazi -> Foolad (reuses an existing shared string)
$ws.Range("S38").Value = "فولاد"
$ws.Range("T38").Value = "Foolad"

# Row 39: Esteghlal Khuzestan -> Sanat Naft Abadan (brand-new strings)
$ws.Range("T39").Value = "Sanat Naft Abadan"
$ws.Range("S39").Value = "صنعت نفت آبادان"

# Row 40: Gostaresh Foolad -> Nassagi Mazandaran (brand-new strings)
$ws.Range("T40").Value = "Nassagi Mazandaran"
$ws.Range("S40").Value = "نساجی مازندران"

# Row 41: Malavan Anzali -> Paykan (brand-new strings)
$ws.Range("T41").Value = "Paykan"
$ws.Range("S41").Value = "پیکان"

# Rows 42-45: English names entered first (brand-new strings) ...
$ws.Range("T42").Value = "Pars Jonobi Jam"
$ws.Range("T43").Value = "Mashin Saazi Tabriz"
$ws.Range("T44").Value = "Naft Masjed Soleiman"
$ws.Range("T45").Value = "Sepid Rood Rasht"

# ... then the Farsi names entered afterwards (brand-new strings)
$ws.Range("S42").Value = "پارس  جنوبی جم"
$ws.Range("S43").Value = "ماشین سازی تبریز"
$ws.Range("S44").Value = "نفت مسجد سلیمان"
$ws.Range("S45").Value = "سپید رود رشت"

# Row 46: Aboumoslem Khorasan -> Esteghlal Khuzestan (reuses an existing shared string)
$ws.Range("S46").Value = "استقلال خوزستان"
$ws.Range("T46").Value = "Esteghlal Khuzestan"

# Update the selection to match the saved view state.
$ws.Range("S45").Select()
